# Re-home the UML "association (partial)" diagram: drop the empty title
# placeholder left over from an earlier layout pass and slide the diagram
# shapes up/right into the freed space (they were anchored to the bottom
# half of the slide while the unused title box sat on top).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# 1. Remove the empty "Title 24" placeholder - it never held any visible
#    text (just an inherited "Click to edit Master title style" prompt)
#    and is no longer needed.
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Title 24") {
        $sh.Delete()
    }
}

# 2. Nudge every remaining diagram shape by the same offset
#    (+9.84378pt horizontally, -119.86087pt vertically == +125016 / -1522233 EMU)
#    so the diagram recenters in the space freed up by the deleted title.
$targets = @{
  "Rectangle 9"        = @(153.84378051757812, 211.05780029296875)
  "Folded Corner 11"   = @(87.84378051757812,  168.13897705078125)
  "TextBox 12"         = @(87.84378051757812,  168.13914489746094)
  "Rectangle 13"       = @(483.8437805175781,  210.13906860351562)
  "Elbow Connector 15" = @(261.8437805175781,  224.22032165527344)
  "TextBox 23"         = @(357.8437805175781,  192.13906860351562)
  "TextBox 34"         = @(261.8437805175781,  228.13906860351562)
  "Line Callout 1 1"   = @(279.8437805175781,  288.13916015625)
}

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($targets.ContainsKey($sh.Name)) {
        $t = $targets[$sh.Name]
        $sh.Left = $t[0]
        $sh.Top  = $t[1]
    }
}

# 3. Presentation-level bookkeeping: slide numbering restarts at 1 instead
#    of being offset to print as "slide 11" (an artifact of the deck this
#    diagram used to live in).
$p.PageSetup.FirstSlideNumber = 1
